$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove now-unused columns I:N entirely (shrinks used range to A:H)
$ws.Range("I1:N1").EntireColumn.Delete()

# 2. Rows that should only span through column E: drop their F:H cells entirely
$ws.Range("F1:H1").Clear()
$ws.Range("F2:H2").Clear()
$ws.Range("F6:H6").Clear()
$ws.Range("F7:H7").Clear()
$ws.Range("F8:H8").Clear()

# 3. Row 4: update the year headers (2008-2011 -> 2017-2020), keep existing style
$ws.Cells.Item(4,4).Value = 2017
$ws.Cells.Item(4,5).Value = 2018
$ws.Cells.Item(4,6).Value = 2019
$ws.Cells.Item(4,7).Value = 2020

# 4. Row 5: update the data values (2017-2021 data)
$ws.Cells.Item(5,4).Value = 0.11
$ws.Cells.Item(5,5).Value = 0.1
$ws.Cells.Item(5,6).Value = 0.09
$ws.Cells.Item(5,7).Value = 0.09
$ws.Cells.Item(5,8).Value = 0.08

# 5. New strings, written in the exact order they should be appended to sharedStrings.xml
$ws.Cells.Item(6,2).Value = "*предварительные данные"
$ws.Cells.Item(4,8).Value = "2021*"
$ws.Cells.Item(6,1).Value = "*алдын алаа маалыматтар"
$ws.Cells.Item(6,3).Value = "*preliminary data"
$ws.Cells.Item(1,1).Value = "9.5.1 ИДП га болгон тажрыйбалык-конструктордук жумуштун жана илимий изилдөөнүн чыгымдарынын үлүшү"
$ws.Cells.Item(1,2).Value = "9.5.1  Доля расходов на научно-исследовательские и опытно-конструкторские работы в ВВП"
$ws.Cells.Item(1,3).Value = "9.5.1 Research and development expenditure as a proportion of GDP"

# 6. New style for H4: right-aligned bold 9pt Times New Roman, keeps the existing medium-bottom border
$h4 = $ws.Cells.Item(4,8)
$h4.HorizontalAlignment = -4152
$h4.VerticalAlignment = -4108
$h4.Font.Bold = $true
$h4.Font.Size = 9
$h4.Font.Name = "Times New Roman"

# 7. Row heights
$ws.Rows.Item(1).RowHeight = 43.5
